$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.120.24"
$ws.Range("E2").Value = "  +5.92%  "

$ws.Range("D3").Value = "3.721.74"
$ws.Range("E3").Value = "  +19.69%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "619.63"
$ws.Range("E5").Value = "  +7.89%  "

$ws.Range("D6").Value = "182.89"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("D7").Value = "3.719.14"
$ws.Range("E7").Value = "  +19.66%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("E10").Value = "  +8.12%  "

$ws.Range("E11").Value = "  +3.64%  "

$ws.Range("E12").Value = "  +7.20%  "

$ws.Range("D13").Value = "40.73"
$ws.Range("E13").Value = "  +12.85%  "

$ws.Range("E14").Value = "  +6.33%  "

$ws.Range("D15").Value = "4.342.27"

$ws.Range("D16").Value = "3.719.13"
$ws.Range("E16").Value = "  +19.76%  "

$ws.Range("D17").Value = "71.141.65"
$ws.Range("E17").Value = "  +6.04%  "

$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("E19").Value = "  +6.96%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "520.24"
$ws.Range("E20").Value = "  +5.10%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "16.94"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").Value = "9.32"
$ws.Range("E22").Value = "  +19.76%  "

$ws.Range("E23").Value = "  +8.51%  "

$ws.Range("E24").Value = "  +12.80%  "

$ws.Range("D25").Value = "88.82"
$ws.Range("E25").Value = "  +6.14%  "

$ws.Range("E26").Value = "  +7.68%  "

$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  +11.18%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "2.55"
$ws.Range("E29").Value = "  +10.29%  "

$ws.Range("D30").Value = "8.23"
$ws.Range("E30").Value = "  +3.99%  "

$ws.Range("E31").Value = "  +11.55%  "

$ws.Range("D32").Value = "32.00"
$ws.Range("E32").Value = "  +13.60%  "

$ws.Range("D33").Value = "0.0000112"
$ws.Range("E33").Value = "  +18.60%  "

$ws.Range("E34").Value = "  +4.55%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "6.16"
$ws.Range("E36").Value = "  +10.23%  "

$ws.Range("E37").Value = "  +10.37%  "

$ws.Range("E38").Value = "  +12.45%  "

$ws.Range("E39").Value = "  +10.95%  "

$ws.Range("D40").Value = "0.134"
$ws.Range("E40").Value = "  +8.99%  "

$ws.Range("D41").Value = "51.63"
$ws.Range("E41").Value = "  +5.18%  "

$ws.Range("D42").Value = "436.56"
$ws.Range("E42").Value = "  +17.25%  "

$ws.Range("D43").Value = "44.97"
$ws.Range("E43").Value = "  -5.59%  "

$ws.Range("D44").Value = "3.152.74"
$ws.Range("E44").Value = "  +12.85%  "

$ws.Range("E45").Value = "  +6.87%  "

$ws.Range("E46").Value = "  +5.55%  "

$ws.Range("D47").Value = "0.0369"
$ws.Range("E47").Value = "  +7.17%  "

$ws.Range("D48").Value = "28.28"
$ws.Range("E48").Value = "  +11.26%  "

$ws.Range("D49").Value = "140.78"
$ws.Range("E49").Value = "  +3.87%  "

$ws.Range("D51").Value = "2.49"
$ws.Range("E51").Value = "  +8.81%  "
